$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "pointofsale1"
$ws.Range("B2").Select()
